$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 1.02
$ws.Cells.Item(2, 3).Value2 = 1.034359590016015
$ws.Cells.Item(2, 4).Value2 = 1.036500297272966
$ws.Cells.Item(2, 5).Value2 = 1.04757400339865
$ws.Cells.Item(2, 6).Value2 = 1.054394617497509
$ws.Cells.Item(2, 9).Value2 = 1.036201307448936
$ws.Cells.Item(2, 10).Value2 = 1.039479264586259
$ws.Cells.Item(2, 11).Value2 = 1.039293861062305
$ws.Cells.Item(2, 12).Value2 = 1.050336261232564
$ws.Cells.Item(2, 13).Value2 = 1.057137947792307
$ws.Cells.Item(2, 14).Value2 = 1.017115050988997

$ws.Cells.Item(3, 2).Value2 = 1.02
$ws.Cells.Item(3, 3).Value2 = 1.03519842061659
$ws.Cells.Item(3, 4).Value2 = 1.037104001315391
$ws.Cells.Item(3, 5).Value2 = 1.048432563529527
$ws.Cells.Item(3, 6).Value2 = 1.05539456535037
$ws.Cells.Item(3, 9).Value2 = 1.036368744643871
$ws.Cells.Item(3, 10).Value2 = 1.039961753859895
$ws.Cells.Item(3, 11).Value2 = 1.039707725807436
$ws.Cells.Item(3, 12).Value2 = 1.051006571847872
$ws.Cells.Item(3, 13).Value2 = 1.057950652928007
$ws.Cells.Item(3, 14).Value2 = 1.017276621537165

$ws.Cells.Item(4, 2).Value2 = 1.02
$ws.Cells.Item(4, 3).Value2 = 1.035741779016724
$ws.Cells.Item(4, 4).Value2 = 1.037495070792547
$ws.Cells.Item(4, 5).Value2 = 1.048989106564996
$ws.Cells.Item(4, 6).Value2 = 1.056042925020959
$ws.Cells.Item(4, 9).Value2 = 1.036476153298078
$ws.Cells.Item(4, 10).Value2 = 1.040273873584554
$ws.Cells.Item(4, 11).Value2 = 1.039975260480038
$ws.Cells.Item(4, 12).Value2 = 1.051440660930273
$ws.Cells.Item(4, 13).Value2 = 1.058477237984291
$ws.Cells.Item(4, 14).Value2 = 1.017381098084007

$ws.Cells.Item(5, 2).Value2 = 1.02
$ws.Cells.Item(5, 3).Value2 = 1.035970343994206
$ws.Cells.Item(5, 4).Value2 = 1.037659578333906
$ws.Cells.Item(5, 5).Value2 = 1.049223314327696
$ws.Cells.Item(5, 6).Value2 = 1.056315811352124
$ws.Cells.Item(5, 9).Value2 = 1.036521083575166
$ws.Cells.Item(5, 10).Value2 = 1.040405067690186
$ws.Cells.Item(5, 11).Value2 = 1.040087667858091
$ws.Cells.Item(5, 12).Value2 = 1.051623235177799
$ws.Cells.Item(5, 13).Value2 = 1.05869878294345
$ws.Cells.Item(5, 14).Value2 = 1.017425002742504

$ws.Cells.Item(6, 2).Value2 = 1.02
$ws.Cells.Item(6, 3).Value2 = 1.036008729056497
$ws.Cells.Item(6, 4).Value2 = 1.037687205803936
$ws.Cells.Item(6, 5).Value2 = 1.049262652695158
$ws.Cells.Item(6, 6).Value2 = 1.056361648645979
$ws.Cells.Item(6, 9).Value2 = 1.036528614394044
$ws.Cells.Item(6, 10).Value2 = 1.040427094491889
$ws.Cells.Item(6, 11).Value2 = 1.040106537754375
$ws.Cells.Item(6, 12).Value2 = 1.051653895035933
$ws.Cells.Item(6, 13).Value2 = 1.058735991160529
$ws.Cells.Item(6, 14).Value2 = 1.017432373502014

$ws.Cells.Item(7, 2).Value2 = 1.02
$ws.Cells.Item(7, 3).Value2 = 1.035744832576261
$ws.Cells.Item(7, 4).Value2 = 1.03749726855189
$ws.Cells.Item(7, 5).Value2 = 1.048992235131145
$ws.Cells.Item(7, 6).Value2 = 1.056046570102164
$ws.Cells.Item(7, 9).Value2 = 1.036476754540842
$ws.Cells.Item(7, 10).Value2 = 1.040275626690652
$ws.Cells.Item(7, 11).Value2 = 1.039976762726582
$ws.Cells.Item(7, 12).Value2 = 1.051443100171124
$ws.Cells.Item(7, 13).Value2 = 1.058480197618254
$ws.Cells.Item(7, 14).Value2 = 1.017381684808644

$ws.Cells.Item(8, 2).Value2 = 1.02
$ws.Cells.Item(8, 3).Value2 = 1.034642955982916
$ws.Cells.Item(8, 4).Value2 = 1.036704231329562
$ws.Cells.Item(8, 5).Value2 = 1.047863950331301
$ws.Cells.Item(8, 6).Value2 = 1.054732278891768
$ws.Cells.Item(8, 9).Value2 = 1.036258086733445
$ws.Cells.Item(8, 10).Value2 = 1.039642340620063
$ws.Cells.Item(8, 11).Value2 = 1.039433782476148
$ws.Cells.Item(8, 12).Value2 = 1.050562721860771
$ws.Cells.Item(8, 13).Value2 = 1.057412457289018
$ws.Cells.Item(8, 14).Value2 = 1.017169668832363

$ws.Cells.Item(9, 2).Value2 = 1.02
$ws.Cells.Item(9, 3).Value2 = 1.032705808217047
$ws.Cells.Item(9, 4).Value2 = 1.03531018361256
$ws.Cells.Item(9, 5).Value2 = 1.045883479600214
$ws.Cells.Item(9, 6).Value2 = 1.052426562768446
$ws.Cells.Item(9, 9).Value2 = 1.035865637248284
$ws.Cells.Item(9, 10).Value2 = 1.038525823323236
$ws.Cells.Item(9, 11).Value2 = 1.038475013303251
$ws.Cells.Item(9, 12).Value2 = 1.049014149728761
$ws.Cells.Item(9, 13).Value2 = 1.055536469928454
$ws.Cells.Item(9, 14).Value2 = 1.016795549295219

$ws.Cells.Item(10, 2).Value2 = 1.02
$ws.Cells.Item(10, 3).Value2 = 1.031417490016368
$ws.Cells.Item(10, 4).Value2 = 1.034383191677419
$ws.Cells.Item(10, 5).Value2 = 1.044568449616703
$ws.Cells.Item(10, 6).Value2 = 1.050896397030114
$ws.Cells.Item(10, 9).Value2 = 1.035599248355794
$ws.Cells.Item(10, 10).Value2 = 1.03778115168146
$ws.Cells.Item(10, 11).Value2 = 1.037834576880744
$ws.Cells.Item(10, 12).Value2 = 1.047983710897435
$ws.Cells.Item(10, 13).Value2 = 1.054289603395053
$ws.Cells.Item(10, 14).Value2 = 1.016545811891743

$ws.Cells.Item(11, 2).Value2 = 1.02
$ws.Cells.Item(11, 3).Value2 = 1.030860391050112
$ws.Cells.Item(11, 4).Value2 = 1.033982377355649
$ws.Cells.Item(11, 5).Value2 = 1.044000301249695
$ws.Cells.Item(11, 6).Value2 = 1.050235494172823
$ws.Cells.Item(11, 9).Value2 = 1.035482778796645
$ws.Cells.Item(11, 10).Value2 = 1.037458637309812
$ws.Cells.Item(11, 11).Value2 = 1.037556976725011
$ws.Cells.Item(11, 12).Value2 = 1.047537997553983
$ws.Cells.Item(11, 13).Value2 = 1.053750614399124
$ws.Cells.Item(11, 14).Value2 = 1.016437601253561

$ws.Cells.Item(12, 2).Value2 = 1.019999999999999
$ws.Cells.Item(12, 3).Value2 = 1.03065357404293
$ws.Cells.Item(12, 4).Value2 = 1.033833585764578
$ws.Cells.Item(12, 5).Value2 = 1.043789457877117
$ws.Cells.Item(12, 6).Value2 = 1.049990257829599
$ws.Cells.Item(12, 9).Value2 = 1.035439348950003
$ws.Cells.Item(12, 10).Value2 = 1.037338832233158
$ws.Cells.Item(12, 11).Value2 = 1.037453821676222
$ws.Cells.Item(12, 12).Value2 = 1.04737251229096
$ws.Cells.Item(12, 13).Value2 = 1.053550548463677
$ws.Cells.Item(12, 14).Value2 = 1.016397396511134

$ws.Cells.Item(13, 2).Value2 = 1.02
$ws.Cells.Item(13, 3).Value2 = 1.030697931794526
$ws.Cells.Item(13, 4).Value2 = 1.033865498025707
$ws.Cells.Item(13, 5).Value2 = 1.043834675770542
$ws.Cells.Item(13, 6).Value2 = 1.050042850409228
$ws.Cells.Item(13, 9).Value2 = 1.035448672385991
$ws.Cells.Item(13, 10).Value2 = 1.037364531212392
$ws.Cells.Item(13, 11).Value2 = 1.037475950666517
$ws.Cells.Item(13, 12).Value2 = 1.047408006137782
$ws.Cells.Item(13, 13).Value2 = 1.053593456991117
$ws.Cells.Item(13, 14).Value2 = 1.016406021032834

$ws.Cells.Item(14, 2).Value2 = 1.02
$ws.Cells.Item(14, 3).Value2 = 1.030843293150443
$ws.Cells.Item(14, 4).Value2 = 1.033970076377384
$ws.Cells.Item(14, 5).Value2 = 1.043982868939628
$ws.Cells.Item(14, 6).Value2 = 1.050215217701622
$ws.Cells.Item(14, 9).Value2 = 1.035479192292849
$ws.Cells.Item(14, 10).Value2 = 1.037448734359379
$ws.Cells.Item(14, 11).Value2 = 1.037548450746611
$ws.Cells.Item(14, 12).Value2 = 1.047524317003155
$ws.Cells.Item(14, 13).Value2 = 1.053734074050936
$ws.Cells.Item(14, 14).Value2 = 1.016434278128595

$ws.Cells.Item(15, 2).Value2 = 1.02
$ws.Cells.Item(15, 3).Value2 = 1.030932870276582
$ws.Cells.Item(15, 4).Value2 = 1.034034522353856
$ws.Cells.Item(15, 5).Value2 = 1.044074201164447
$ws.Cells.Item(15, 6).Value2 = 1.050321452374965
$ws.Cells.Item(15, 9).Value2 = 1.035497974389164
$ws.Cells.Item(15, 10).Value2 = 1.037500613554621
$ws.Cells.Item(15, 11).Value2 = 1.037593114908125
$ws.Cells.Item(15, 12).Value2 = 1.047595989616859
$ws.Cells.Item(15, 13).Value2 = 1.053820731269811
$ws.Cells.Item(15, 14).Value2 = 1.016451686880234

$ws.Cells.Item(16, 2).Value2 = 1.02
$ws.Cells.Item(16, 3).Value2 = 1.031454478795343
$ws.Cells.Item(16, 4).Value2 = 1.034409804750952
$ws.Cells.Item(16, 5).Value2 = 1.044606182617601
$ws.Cells.Item(16, 6).Value2 = 1.050940294281358
$ws.Cells.Item(16, 9).Value2 = 1.035606954482269
$ws.Cells.Item(16, 10).Value2 = 1.037802554590799
$ws.Cells.Item(16, 11).Value2 = 1.037852994338657
$ws.Cells.Item(16, 12).Value2 = 1.048013301517732
$ws.Cells.Item(16, 13).Value2 = 1.054325393670975
$ws.Cells.Item(16, 14).Value2 = 1.016552991980379

$ws.Cells.Item(17, 2).Value2 = 1.02
$ws.Cells.Item(17, 3).Value2 = 1.031781872228497
$ws.Cells.Item(17, 4).Value2 = 1.034645365752701
$ws.Cells.Item(17, 5).Value2 = 1.044940221272697
$ws.Cells.Item(17, 6).Value2 = 1.051328925467085
$ws.Cells.Item(17, 9).Value2 = 1.03567501501383
$ws.Cells.Item(17, 10).Value2 = 1.037991937184059
$ws.Cells.Item(17, 11).Value2 = 1.038015933804104
$ws.Cells.Item(17, 12).Value2 = 1.048275197986132
$ws.Cells.Item(17, 13).Value2 = 1.054642200400956
$ws.Cells.Item(17, 14).Value2 = 1.016616518833177

$ws.Cells.Item(18, 2).Value2 = 1.02
$ws.Cells.Item(18, 3).Value2 = 1.031972907670856
$ws.Cells.Item(18, 4).Value2 = 1.034782820282784
$ws.Cells.Item(18, 5).Value2 = 1.045135182627852
$ws.Cells.Item(18, 6).Value2 = 1.051555768386125
$ws.Cells.Item(18, 9).Value2 = 1.03571460532931
$ws.Cells.Item(18, 10).Value2 = 1.03810239434113
$ws.Cells.Item(18, 11).Value2 = 1.038110945875008
$ws.Cells.Item(18, 12).Value2 = 1.048428003355272
$ws.Cells.Item(18, 13).Value2 = 1.054827076343694
$ws.Cells.Item(18, 14).Value2 = 1.016653565928953

$ws.Cells.Item(19, 2).Value2 = 1.02
$ws.Cells.Item(19, 3).Value2 = 1.032038058053943
$ws.Cells.Item(19, 4).Value2 = 1.034829698136634
$ws.Cells.Item(19, 5).Value2 = 1.045201680097863
$ws.Cells.Item(19, 6).Value2 = 1.051633143204487
$ws.Cells.Item(19, 9).Value2 = 1.035728086225606
$ws.Cells.Item(19, 10).Value2 = 1.038140056239349
$ws.Cells.Item(19, 11).Value2 = 1.038143337792728
$ws.Cells.Item(19, 12).Value2 = 1.048480113728077
$ws.Cells.Item(19, 13).Value2 = 1.054890129132852
$ws.Cells.Item(19, 14).Value2 = 1.016666196823641

$ws.Cells.Item(20, 2).Value2 = 1.02
$ws.Cells.Item(20, 3).Value2 = 1.031746738513047
$ws.Cells.Item(20, 4).Value2 = 1.034620086504645
$ws.Cells.Item(20, 5).Value2 = 1.044904369428347
$ws.Cells.Item(20, 6).Value2 = 1.051287212361797
$ws.Cells.Item(20, 9).Value2 = 1.035667723954376
$ws.Cells.Item(20, 10).Value2 = 1.037971618888819
$ws.Cells.Item(20, 11).Value2 = 1.037998454805024
$ws.Cells.Item(20, 12).Value2 = 1.048247094248048
$ws.Cells.Item(20, 13).Value2 = 1.05460820093054
$ws.Cells.Item(20, 14).Value2 = 1.01660970372845

$ws.Cells.Item(21, 2).Value2 = 1.02
$ws.Cells.Item(21, 3).Value2 = 1.030800484685559
$ws.Cells.Item(21, 4).Value2 = 1.033939278197641
$ws.Cells.Item(21, 5).Value2 = 1.043939224425014
$ws.Cells.Item(21, 6).Value2 = 1.050164452849232
$ws.Cells.Item(21, 9).Value2 = 1.035470209569975
$ws.Cells.Item(21, 10).Value2 = 1.037423938867306
$ws.Cells.Item(21, 11).Value2 = 1.037527102430183
$ws.Cells.Item(21, 12).Value2 = 1.047490064339289
$ws.Cells.Item(21, 13).Value2 = 1.053692661994227
$ws.Cells.Item(21, 14).Value2 = 1.016425957404412

$ws.Cells.Item(22, 2).Value2 = 1.02
$ws.Cells.Item(22, 3).Value2 = 1.030206199696972
$ws.Cells.Item(22, 4).Value2 = 1.033511740971945
$ws.Cells.Item(22, 5).Value2 = 1.043333512458791
$ws.Cells.Item(22, 6).Value2 = 1.049459990539675
$ws.Cells.Item(22, 9).Value2 = 1.035345053517825
$ws.Cells.Item(22, 10).Value2 = 1.037079539890281
$ws.Cells.Item(22, 11).Value2 = 1.037230502093833
$ws.Cells.Item(22, 12).Value2 = 1.047014509843629
$ws.Cells.Item(22, 13).Value2 = 1.053117828827147
$ws.Cells.Item(22, 14).Value2 = 1.016310368298381

$ws.Cells.Item(23, 2).Value2 = 1.02
$ws.Cells.Item(23, 3).Value2 = 1.030521178136329
$ws.Cells.Item(23, 4).Value2 = 1.033738337271501
$ws.Cells.Item(23, 5).Value2 = 1.043654505770839
$ws.Cells.Item(23, 6).Value2 = 1.04983330028401
$ws.Cells.Item(23, 9).Value2 = 1.035411492935886
$ws.Cells.Item(23, 10).Value2 = 1.037262116742029
$ws.Cells.Item(23, 11).Value2 = 1.037387758145207
$ws.Cells.Item(23, 12).Value2 = 1.047266569981429
$ws.Cells.Item(23, 13).Value2 = 1.053422482198376
$ws.Cells.Item(23, 14).Value2 = 1.016371649862169

$ws.Cells.Item(24, 2).Value2 = 1.02
$ws.Cells.Item(24, 3).Value2 = 1.031762613703235
$ws.Cells.Item(24, 4).Value2 = 1.03463150893656
$ws.Cells.Item(24, 5).Value2 = 1.044920568956026
$ws.Cells.Item(24, 6).Value2 = 1.051306060221519
$ws.Cells.Item(24, 9).Value2 = 1.035671018804822
$ws.Cells.Item(24, 10).Value2 = 1.037980799872082
$ws.Cells.Item(24, 11).Value2 = 1.038006352898251
$ws.Cells.Item(24, 12).Value2 = 1.048259792976874
$ws.Cells.Item(24, 13).Value2 = 1.054623563556191
$ws.Cells.Item(24, 14).Value2 = 1.016612783202736

$ws.Cells.Item(25, 2).Value2 = 1.02
$ws.Cells.Item(25, 3).Value2 = 1.0332060650324
$ws.Cells.Item(25, 4).Value2 = 1.035670167177124
$ws.Cells.Item(25, 5).Value2 = 1.046394554715993
$ws.Cells.Item(25, 6).Value2 = 1.053021422886544
$ws.Cells.Item(25, 9).Value2 = 1.035967935834997
$ws.Cells.Item(25, 10).Value2 = 1.038814532345951
$ws.Cells.Item(25, 11).Value2 = 1.038723104544298
$ws.Cells.Item(25, 12).Value2 = 1.049414156379645
$ws.Cells.Item(25, 13).Value2 = 1.056020795824879
$ws.Cells.Item(25, 14).Value2 = 1.016892327199673
